$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the MCQ-only columns (Correct/Choice_1..4) that are no longer needed
# and replace column C header/content with the new "Answer" column.
$ws.Range("D1:G1").EntireColumn.Delete()

# Update header row
$ws.Range("C1").Value = "Answer"

# Row 2
$ws.Range("A2").Value = "OpenQuestion"
$ws.Range("B2").Value = "What is the first step in the Git Workflow described in the devops automation lecture notes?"
$ws.Range("C2").Value = "Clone code from the project’s shared repository hosted on a server, e.g., GitHub ."

# Row 3
$ws.Range("A3").Value = "OpenQuestion"
$ws.Range("B3").Value = "In the Git Workflow example provided, what is the purpose of creating multiple branches for different users working on the project?"
$ws.Range("C3").Value = "Each user can work on their own branch to experiment with new features without affecting others and to maintain flexibility in collaboration ."

# Row 4
$ws.Range("A4").Value = "OpenQuestion"
$ws.Range("B4").Value = "How is a critical issue handled in the Git Workflow example when working on a project?"
$ws.Range("C4").Value = "Switch to the production branch, create a new branch to add the fix, test the fix, merge the fix branch, and push to production ."

# Row 5
$ws.Range("A5").Value = "OpenQuestion"
$ws.Range("B5").Value = "What is the purpose of a three-way merge in Git, as described in the lecture notes?"
$ws.Range("C5").Value = "A three-way merge is used to merge branches where a common ancestor is considered along with the two branches to be merged, creating a new commit from the merge ."

# Row 6
$ws.Range("A6").Value = "OpenQuestion"
$ws.Range("B6").Value = "What is the significance of Continuous Integration (CI) implementation in the Git Workflow process?"
$ws.Range("C6").Value = "CI implementation involves monitoring and reacting to new commits, triggering the build/test process automatically on every push to the main repository, and ensuring integration of changes daily ."
